$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "p_adj" column header in H1
$ws.Cells.Item(1, 8).Value = "p_adj"

# p_adj values for rows 2-13 (column H), written as plain decimals
# (shortest round-tripping decimal representation of each double)
$ws.Cells.Item(2, 8).Value = 0.000000652020275104803
$ws.Cells.Item(3, 8).Value = 0.121283397555078
$ws.Cells.Item(4, 8).Value = 0.00000474114302284612
$ws.Cells.Item(5, 8).Value = 0.0000065436111745143
$ws.Cells.Item(6, 8).Value = 0.0000000000173394631985957
$ws.Cells.Item(7, 8).Value = 0.64635519553949
$ws.Cells.Item(8, 8).Value = 0.459279706662271
$ws.Cells.Item(9, 8).Value = 0.0000254993498423062
$ws.Cells.Item(10, 8).Value = 0.0000000000173394631985957
$ws.Cells.Item(11, 8).Value = 0.64635519553949
$ws.Cells.Item(12, 8).Value = 0.0800683006875426
$ws.Cells.Item(13, 8).Value = 0.40742527308503

# Every row that used to be highlighted in column E loses that highlight
# (the significance marker now lives on the adjusted p-value in column H)
$oldHighlightRows = @(2, 4, 5, 6, 9, 10, 12)
foreach ($r in $oldHighlightRows) {
    $ws.Cells.Item($r, 5).ClearFormats()
}

# Rows that remain significant after adjustment get the yellow fill in column H
$newHighlightRows = @(2, 4, 5, 6, 9, 10)
foreach ($r in $newHighlightRows) {
    $ws.Cells.Item($r, 8).Interior.Color = 65535
}

# Update the view: zoom to 100% and move the active selection to H19
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("H19").Select()
